$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04763786555579896
$ws.Range("C2").Value = 0.002777888934908601
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 9.516137877023681

$ws.Range("B3").Value = 0.01514828764759746
$ws.Range("C3").Value = 0.00007097389502863649
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 9.480941384075599

$ws.Range("B4").Value = 0.003994804209775715
$ws.Range("C4").Value = 0.002777888934908601
$ws.Range("D4").Value = 0.8054896365839992
$ws.Range("E4").Value = 645.3272768299601
$ws.Range("G4").Value = 646.1395391596888

$ws.Range("B5").Value = 0.127881588408715
$ws.Range("C5").Value = 1.667794583268128
$ws.Range("D5").Value = 0.8054896365839992
$ws.Range("E5").Value = 8.660232485948974
$ws.Range("G5").Value = 11.26139829420982

$ws.Range("B6").Value = 3.230985683306322
$ws.Range("C6").Value = 1.667794583268128
$ws.Range("D6").Value = 0.1575252929769615
$ws.Range("E6").Value = 0.496779210170732
$ws.Range("G6").Value = 5.553084769722144

$ws.Range("B7").Value = 3.230985683306322
$ws.Range("C7").Value = 1.667794583268128
$ws.Range("D7").Value = 0.1575252929769615
$ws.Range("E7").Value = 0.496779210170732
$ws.Range("G7").Value = 5.553084769722144

$ws.Range("B8").Value = 1.459612070389937
$ws.Range("C8").Value = 1.667794583268128
$ws.Range("D8").Value = 0.1575252929769615
$ws.Range("E8").Value = 0.496779210170732
$ws.Range("G8").Value = 3.781711156805759
